$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.031554796368141
$ws.Range("D2").Value = 1.03669790318922
$ws.Range("E2").Value = 1.031111575187526
$ws.Range("F2").Value = 1.042140145816711
$ws.Range("I2").Value = 1.038605971110588
$ws.Range("J2").Value = 1.036689879752767
$ws.Range("K2").Value = 1.039490901972486
$ws.Range("L2").Value = 1.03392063742663
$ws.Range("M2").Value = 1.044917675235453

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.032460036712933
$ws.Range("D3").Value = 1.037400536042004
$ws.Range("E3").Value = 1.03187961037638
$ws.Range("F3").Value = 1.043306765786308
$ws.Range("I3").Value = 1.038879744769842
$ws.Range("J3").Value = 1.037237269270151
$ws.Range("K3").Value = 1.040003473749549
$ws.Range("L3").Value = 1.03449727674933
$ws.Range("M3").Value = 1.045894133353827

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.033046064061433
$ws.Range("D4").Value = 1.037855333895837
$ws.Range("E4").Value = 1.032377192438332
$ws.Range("F4").Value = 1.044062255530107
$ws.Range("I4").Value = 1.039055692806187
$ws.Range("J4").Value = 1.037591117351046
$ws.Range("K4").Value = 1.040334615681425
$ws.Range("L4").Value = 1.034870353501824
$ws.Range("M4").Value = 1.046526001143164

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.033292495258369
$ws.Range("D5").Value = 1.038046565017826
$ws.Range("E5").Value = 1.032586521217617
$ws.Range("F5").Value = 1.044380009051141
$ws.Range("I5").Value = 1.039129373161291
$ws.Range("J5").Value = 1.037739790641141
$ws.Range("K5").Value = 1.04047370094398
$ws.Range("L5").Value = 1.035027182730889
$ws.Range("M5").Value = 1.046791646328927

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.03333387594523
$ws.Range("D6").Value = 1.038078675530017
$ws.Range("E6").Value = 1.032621676924615
$ws.Range("F6").Value = 1.044433369803565
$ws.Range("I6").Value = 1.039141727503068
$ws.Range("J6").Value = 1.0377647485668
$ws.Range("K6").Value = 1.040497046512009
$ws.Range("L6").Value = 1.035053514305365
$ws.Range("M6").Value = 1.04683624978616

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.033049356631621
$ws.Range("D7").Value = 1.037857889002933
$ws.Range("E7").Value = 1.032379988930405
$ws.Range("F7").Value = 1.044066500796046
$ws.Range("I7").Value = 1.039056678459721
$ws.Range("J7").Value = 1.037593104263583
$ws.Range("K7").Value = 1.040336474644788
$ws.Range("L7").Value = 1.034872449110846
$ws.Range("M7").Value = 1.046529550676067

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.031860668856281
$ws.Range("D8").Value = 1.036935330015097
$ws.Range("E8").Value = 1.031371008837343
$ws.Range("F8").Value = 1.04253428443585
$ws.Range("I8").Value = 1.038698742758257
$ws.Range("J8").Value = 1.03687494451887
$ws.Range("K8").Value = 1.039664236486102
$ws.Range("L8").Value = 1.034115524681801
$ws.Range("M8").Value = 1.045247666650594

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.029768202768907
$ws.Range("D9").Value = 1.03531084447317
$ws.Range("E9").Value = 1.029597797391679
$ws.Range("F9").Value = 1.039838987484017
$ws.Range("I9").Value = 1.038058830847039
$ws.Range("J9").Value = 1.035606811775506
$ws.Range("K9").Value = 1.038475672004901
$ws.Range("L9").Value = 1.032781402400711
$ws.Range("M9").Value = 1.042989092638506

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.028374722259588
$ws.Range("D10").Value = 1.034228717871831
$ws.Range("E10").Value = 1.028418915128055
$ws.Range("F10").Value = 1.038045259478133
$ws.Range("I10").Value = 1.037626076192852
$ws.Range("J10").Value = 1.034759656853621
$ws.Range("K10").Value = 1.037680656988666
$ws.Range("L10").Value = 1.031891820509401
$ws.Range("M10").Value = 1.041483566461935

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.027771694656534
$ws.Range("D11").Value = 1.033760364356198
$ws.Range("E11").Value = 1.027909233831073
$ws.Range("F11").Value = 1.037269298606128
$ws.Range("I11").Value = 1.037437236333834
$ws.Range("J11").Value = 1.034392427367842
$ws.Range("K11").Value = 1.037335790519506
$ws.Range("L11").Value = 1.031506592652149
$ws.Range("M11").Value = 1.040831703029612

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.027547757965036
$ws.Range("D12").Value = 1.033586430349573
$ws.Range("E12").Value = 1.027720034233279
$ws.Range("F12").Value = 1.036981182460219
$ws.Range("I12").Value = 1.037366874725869
$ws.Range("J12").Value = 1.034255961699325
$ws.Range("K12").Value = 1.037207599388003
$ws.Range("L12").Value = 1.03136349772604
$ws.Range("M12").Value = 1.040589577811385

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.027595790662621
$ws.Range("D13").Value = 1.03362373825757
$ws.Range("E13").Value = 1.027760612797778
$ws.Range("F13").Value = 1.037042979364403
$ws.Range("I13").Value = 1.037381977391491
$ws.Range("J13").Value = 1.034285236765603
$ws.Range("K13").Value = 1.037235100988967
$ws.Range("L13").Value = 1.031394192248058
$ws.Range("M13").Value = 1.040641514226667

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.027753182861516
$ws.Range("D14").Value = 1.033745986230276
$ws.Range("E14").Value = 1.027893592104508
$ws.Range("F14").Value = 1.037245480587156
$ws.Range("I14").Value = 1.037431424667702
$ws.Range("J14").Value = 1.034381148298992
$ws.Range("K14").Value = 1.037325196080641
$ws.Range("L14").Value = 1.031494764460821
$ws.Range("M14").Value = 1.040811688763622

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.027850164648004
$ws.Range("D15").Value = 1.033821311692773
$ws.Range("E15").Value = 1.027975540811951
$ws.Range("F15").Value = 1.037370262880578
$ws.Range("I15").Value = 1.037461861881516
$ws.Range("J15").Value = 1.034440234592053
$ws.Range("K15").Value = 1.037380694422283
$ws.Range("L15").Value = 1.031556729795156
$ws.Range("M15").Value = 1.040916539691432

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.028414750833937
$ws.Range("D16").Value = 1.034259805566317
$ws.Range("E16").Value = 1.028452757583249
$ws.Range("F16").Value = 1.038096772919566
$ws.Range("I16").Value = 1.037638578267912
$ws.Range("J16").Value = 1.034784020173352
$ws.Range("K16").Value = 1.037703531646486
$ws.Range("L16").Value = 1.031917386187237
$ws.Range("M16").Value = 1.041526829317604

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.028768997365914
$ws.Range("D17").Value = 1.034534919395727
$ws.Range("E17").Value = 1.028752313378788
$ws.Range("F17").Value = 1.038552690174112
$ws.Range("I17").Value = 1.037749038734095
$ws.Range("J17").Value = 1.034999559619519
$ws.Range("K17").Value = 1.037905873511057
$ws.Range("L17").Value = 1.0321436081462
$ws.Range("M17").Value = 1.041909658389057

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.028975657930608
$ws.Range("D18").Value = 1.034695409253265
$ws.Range("E18").Value = 1.028927114628223
$ws.Range("F18").Value = 1.038818689942444
$ws.Range("I18").Value = 1.037813328114948
$ws.Range("J18").Value = 1.035125240930872
$ws.Range("K18").Value = 1.038023836253867
$ws.Range("L18").Value = 1.032275556346191
$ws.Range("M18").Value = 1.042132959906379

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.029046129681805
$ws.Range("D19").Value = 1.034750135609733
$ws.Range("E19").Value = 1.028986730094346
$ws.Range("F19").Value = 1.038909401028914
$ws.Range("I19").Value = 1.037835225315949
$ws.Range("J19").Value = 1.035168088340225
$ws.Range("K19").Value = 1.038064048327843
$ws.Range("L19").Value = 1.032320546697973
$ws.Range("M19").Value = 1.042209100598012

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.028730986486002
$ws.Range("D20").Value = 1.034505400144978
$ws.Range("E20").Value = 1.028720166067185
$ws.Range("F20").Value = 1.038503767230245
$ws.Range("I20").Value = 1.037737201884548
$ws.Range("J20").Value = 1.034976438326566
$ws.Range("K20").Value = 1.037884170328469
$ws.Range("L20").Value = 1.032119337000297
$ws.Range("M20").Value = 1.041868584060095

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.027706833268484
$ws.Range("D21").Value = 1.033709986319825
$ws.Range("E21").Value = 1.027854429726218
$ws.Range("F21").Value = 1.037185845981668
$ws.Range("I21").Value = 1.037416869693997
$ws.Range("J21").Value = 1.034352906401665
$ws.Range("K21").Value = 1.037298667865646
$ws.Range("L21").Value = 1.031465148558816
$ws.Range("M21").Value = 1.040761576441708

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.027063223810968
$ws.Range("D22").Value = 1.033210071542979
$ws.Range("E22").Value = 1.027310794436373
$ws.Range("F22").Value = 1.036357854495818
$ws.Range("I22").Value = 1.03721420261248
$ws.Range("J22").Value = 1.033960518670941
$ws.Range("K22").Value = 1.036930005170764
$ws.Range("L22").Value = 1.031053810726771
$ws.Range("M22").Value = 1.04006559085349

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.027404383135193
$ws.Range("D23").Value = 1.033475067111603
$ws.Range("E23").Value = 1.027598920232239
$ws.Range("F23").Value = 1.036796728041537
$ws.Range("I23").Value = 1.03732175970827
$ws.Range("J23").Value = 1.03416856363342
$ws.Range("K23").Value = 1.03712549061733
$ws.Range("L23").Value = 1.031271870658958
$ws.Range("M23").Value = 1.040434542763046

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.02874816185972
$ws.Range("D24").Value = 1.034518738560152
$ws.Range("E24").Value = 1.0287346918206
$ws.Range("F24").Value = 1.038525873183321
$ws.Range("I24").Value = 1.037742550880986
$ws.Range("J24").Value = 1.034986885964768
$ws.Range("K24").Value = 1.03789397724743
$ws.Range("L24").Value = 1.032130304096467
$ws.Range("M24").Value = 1.041887143769457

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.030308895160678
$ws.Range("D25").Value = 1.035730665764296
$ws.Range("E25").Value = 1.030055645916074
$ws.Range("F25").Value = 1.040535233574657
$ws.Range("I25").Value = 1.038225348850053
$ws.Range("J25").Value = 1.035934963196349
$ws.Range("K25").Value = 1.038783412961641
$ws.Range("L25").Value = 1.033126338183153
$ws.Range("M25").Value = 1.043572955349557
